# ---------------------------------------------------------------------------
# api_test_cases_multi_flow.xlsx: turn the single "Sheet1" login-test sheet
# into a multi-sheet "restful API framework" workbook:
#   pre_process, after_process, users, products, orders, payments,
#   shipments, Sheet6
# and replace the data on the (renamed) "users" sheet with the new
# case_id/interface/title/method/url/request_data/expected/extract_data
# table.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Insert the two sheets that must land BEFORE the original sheet ----
# Worksheets.Add(Before) puts the new sheet immediately before the sheet
# handed in. Re-fetch Item(1) fresh each time (object handles don't track
# identity across structural edits in this host, but positional Item()
# lookups do).
$null = $wb.Worksheets.Add($wb.Worksheets.Item(1))   # -> Sheet2, Sheet1
$null = $wb.Worksheets.Add($wb.Worksheets.Item(1))   # -> Sheet3, Sheet2, Sheet1

# Original data sheet is now at position 3. Append the five trailing sheets
# right after it (and after each other), one at a time.
for ($k = 1; $k -le 5; $k++) {
    $anchor = $wb.Worksheets.Item(2 + $k)
    $null = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $anchor)
}

# --- 2. Name every tab in final left-to-right order ------------------------
$names = @("pre_process", "after_process", "users", "products", "orders", "payments", "shipments", "Sheet6")
for ($i = 1; $i -le 8; $i++) {
    $wb.Worksheets.Item($i).Name = $names[$i - 1]
}

$preProcess   = $wb.Worksheets.Item(1)
$afterProcess = $wb.Worksheets.Item(2)
$users        = $wb.Worksheets.Item(3)
$products     = $wb.Worksheets.Item(4)
$orders       = $wb.Worksheets.Item(5)
$payments     = $wb.Worksheets.Item(6)
$shipments    = $wb.Worksheets.Item(7)
$sheet6       = $wb.Worksheets.Item(8)

# --- 3. Rewrite the "users" sheet's table -----------------------------------
# header row
$users.Range("A1").Value = "case_id"
$users.Range("B1").Value = "interface"
$users.Range("C1").Value = "title"
$users.Range("D1").Value = "method"
$users.Range("E1").Value = "url"
$users.Range("F1").Value = "request_data"
$users.Range("G1").Value = "expected"
$users.Range("H1").Value = "extract_data"

# row 2 - add user
$users.Range("A2").Value = 1
$users.Range("B2").Value = "users_add"
$users.Range("C2").Value = "user added successfully"
$users.Range("D2").Value = "post"
$users.Range("E2").Value = "/paymall_admin/users/"
$users.Range("F2").Value = '{"username":"#username#","mobile":"#phone#","password":"#password#","email":"a@a.com"}'
$users.Range("G2").Value = '{"id":#id#,"username":"#username#","mobile":"#phone#","email":"a@a.com"}'

# row 3 - login
$users.Range("A3").Value = 2
$users.Range("B3").Value = "login"
$users.Range("C3").Value = "login pass, username, password"
$users.Range("D3").Value = "post"
$users.Range("E3").Value = "/paymall_admin/authorizations/"
$users.Range("F3").Value = '{"username": "#username#", "password": "#password#"}'
$users.Range("G3").Value = '{"refresh":"#refresh_token#", "access":"#access_token#"}'

# row 4 - list current users (no request_data column, F4 stays empty)
$users.Range("A4").Value = 3
$users.Range("B4").Value = "users_list"
$users.Range("C4").Value = "list current users"
$users.Range("D4").Value = "get"
$users.Range("E4").Value = "/paymall_admin/statistical/total_count/"
$users.Range("G4").Value = '{"count":#count#}'

# --- 4. Column widths on "users" (engine quantises ColumnWidth to the
#        nearest 1/6 of a character, so feed in pre-snapped values to land
#        as close as possible to the authored widths) --------------------
$users.Columns.Item(1).ColumnWidth = 7.666666666666667   # -> width 8.5
$users.Columns.Item(2).ColumnWidth = 13.333333333333334  # -> width 14.1666...
$users.Columns.Item(3).ColumnWidth = 29.5                # -> width 30.3333...
$users.Columns.Item(5).ColumnWidth = 36.333333333333336  # -> width 37.1666...
$users.Columns.Item(6).ColumnWidth = 85.5                # -> width 86.3333...
$users.Columns.Item(7).ColumnWidth = 58.0                # -> width 58.8333...
$users.Columns.Item(8).ColumnWidth = 22.166666666666668  # -> width 23

# --- 5. Sheet views / selections --------------------------------------------
$null = $preProcess.Range("G15").Select()
$null = $afterProcess.Range("D13").Select()
$null = $shipments.Range("H16").Select()
$null = $sheet6.Range("H13").Select()

# "users" is the tab that ends up active/selected, scrolled so column D is
# the left-most visible column, with G8 as the active cell.
$null = $users.Range("G8").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 4

$users.Activate()
